$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '198000'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '142880'
$ws.Range("B15").Value = '삼성전자 삼성 SL-J1770FW 정품잉크포함'
$ws.Range("C15").Value = 'https://search.shopping.naver.com/gate.nhn?id=18378943253'
$ws.Range("D15").Value = 'https://shopping-phinf.pstatic.net/main_1837894/18378943253.20190329110354.jpg'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '139000'
$ws.Range("I15").Value = '삼성'
$ws.Range("J15").Value = '삼성전자'
$ws.Range("M15").Value = '복합기'
$ws.Range("N15").Value = '잉크젯복합기'
$ws.Range("B16").Value = '사나고 3D펜 고급형'
$ws.Range("C16").Value = 'https://search.shopping.naver.com/gate.nhn?id=25665893522'
$ws.Range("D16").Value = 'https://shopping-phinf.pstatic.net/main_2566589/25665893522.20210119160253.jpg'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '59000'
$ws.Range("I16").Value = '사나고'
$ws.Range("J16").Value = ''
$ws.Range("M16").Value = '프린터'
$ws.Range("N16").Value = '3D프린터'
$ws.Range("B19").Value = '삼성전자 삼성 SL-M2030W 정품토너포함'
$ws.Range("C19").Value = 'https://search.shopping.naver.com/gate.nhn?id=21379907273'
$ws.Range("D19").Value = 'https://shopping-phinf.pstatic.net/main_2137990/21379907273.20200221180856.jpg'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '134900'
$ws.Range("B20").Value = '삼성전자 삼성 SL-M2843DW 정품토너포함'
$ws.Range("C20").Value = 'https://search.shopping.naver.com/gate.nhn?id=22846051427'
$ws.Range("D20").Value = 'https://shopping-phinf.pstatic.net/main_2284605/22846051427.20200724162901.jpg'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '209000'
$ws.Range("B27").Value = 'HP 2131 가정용 프린터기 무한잉크 복합기 잉크젯 프린트 복사 공기계'
$ws.Range("C27").Value = 'https://search.shopping.naver.com/gate.nhn?id=82335217020'
$ws.Range("D27").Value = 'https://shopping-phinf.pstatic.net/main_8233521/82335217020.jpg'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '49000'
$ws.Range("G27").Value = '잉크시대'
$ws.Range("H27").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I27").Value = 'HP'
$ws.Range("J27").Value = 'HP'
$ws.Range("M27").Value = '복합기'
$ws.Range("N27").Value = '잉크젯복합기'
$ws.Range("B28").Value = '아이티씨 퓨리웨이 휴대용 무선 라벨프린터 D11'
$ws.Range("C28").Value = 'https://search.shopping.naver.com/gate.nhn?id=28177801522'
$ws.Range("D28").Value = 'https://shopping-phinf.pstatic.net/main_2817780/28177801522.20210727183423.jpg'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '39000'
$ws.Range("G28").Value = '네이버'
$ws.Range("H28").Value = '일반 - 가격비교 상품'
$ws.Range("I28").Value = '퓨리웨이'
$ws.Range("J28").Value = '아이티씨'
$ws.Range("M28").Value = '프린터'
$ws.Range("N28").Value = '라벨프린터'
$ws.Range("B29").Value = '엡손 완성형 L3156 정품 무한잉크'
$ws.Range("C29").Value = 'https://search.shopping.naver.com/gate.nhn?id=17767089381'
$ws.Range("D29").Value = 'https://shopping-phinf.pstatic.net/main_1776708/17767089381.20210701153200.jpg'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '199000'
$ws.Range("I29").Value = '엡손'
$ws.Range("J29").Value = '엡손'
$ws.Range("M29").Value = '복합기'
$ws.Range("N29").Value = '잉크젯복합기'
$ws.Range("B30").Value = '디앤에스글로벌 휴대용 멀티 스티커 프린터 프릭커'
$ws.Range("C30").Value = 'https://search.shopping.naver.com/gate.nhn?id=28859379589'
$ws.Range("D30").Value = 'https://shopping-phinf.pstatic.net/main_2885937/28859379589.20210915163615.jpg'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '159000'
$ws.Range("I30").Value = ''
$ws.Range("J30").Value = '디앤에스글로벌'
$ws.Range("M30").Value = '프린터'
$ws.Range("N30").Value = '라벨프린터'
$ws.Range("B62").Value = '브라더 PT-P300BT'
$ws.Range("C62").Value = 'https://search.shopping.naver.com/gate.nhn?id=13419538634'
$ws.Range("D62").Value = 'https://shopping-phinf.pstatic.net/main_1341953/13419538634.20200108103842.jpg'
$ws.Range("E62").NumberFormat = "@"
$ws.Range("E62").Value = '79000'
$ws.Range("I62").Value = '브라더'
$ws.Range("J62").Value = '브라더'
$ws.Range("B63").Value = '엡손 LW-K200KP'
$ws.Range("C63").Value = 'https://search.shopping.naver.com/gate.nhn?id=26208478522'
$ws.Range("D63").Value = 'https://shopping-phinf.pstatic.net/main_2620847/26208478522.20210402090844.jpg'
$ws.Range("E63").NumberFormat = "@"
$ws.Range("E63").Value = '95920'
$ws.Range("I63").Value = '엡손'
$ws.Range("J63").Value = '엡손'
$ws.Range("B70").Value = '카피어랜드 D30S'
$ws.Range("C70").Value = 'https://search.shopping.naver.com/gate.nhn?id=27505417522'
$ws.Range("D70").Value = 'https://shopping-phinf.pstatic.net/main_2750541/27505417522.20211124163427.jpg'
$ws.Range("E70").NumberFormat = "@"
$ws.Range("E70").Value = '38900'
$ws.Range("I70").Value = '카피어랜드'
$ws.Range("J70").Value = '카피어랜드'
$ws.Range("N70").Value = '라벨프린터'
$ws.Range("B71").Value = '캐논 인스픽 iNSPiC PV-123'
$ws.Range("C71").Value = 'https://search.shopping.naver.com/gate.nhn?id=16198366421'
$ws.Range("D71").Value = 'https://shopping-phinf.pstatic.net/main_1619836/16198366421.20190115111910.jpg'
$ws.Range("E71").NumberFormat = "@"
$ws.Range("E71").Value = '129000'
$ws.Range("I71").Value = '캐논'
$ws.Range("J71").Value = '캐논'
$ws.Range("N71").Value = '포토프린터'
$ws.Range("B75").Value = '캐논 G3960 정품 무한잉크'
$ws.Range("C75").Value = 'https://search.shopping.naver.com/gate.nhn?id=25732367524'
$ws.Range("D75").Value = 'https://shopping-phinf.pstatic.net/main_2573236/25732367524.20210125110819.jpg'
$ws.Range("E75").NumberFormat = "@"
$ws.Range("E75").Value = '219000'
$ws.Range("G75").Value = '네이버'
$ws.Range("H75").Value = '일반 - 가격비교 상품'
$ws.Range("I75").Value = '캐논'
$ws.Range("J75").Value = '캐논'
$ws.Range("M75").Value = '복합기'
$ws.Range("N75").Value = '잉크젯복합기'
$ws.Range("B76").Value = '삼성 프린터 컬러 레이저 프린트 가정용 사무실 무선 와이파이 레이져 복합기 팩스 스캔 복사'
$ws.Range("C76").Value = 'https://search.shopping.naver.com/gate.nhn?id=82314641450'
$ws.Range("D76").Value = 'https://shopping-phinf.pstatic.net/main_8231464/82314641450.1.jpg'
$ws.Range("E76").NumberFormat = "@"
$ws.Range("E76").Value = '179000'
$ws.Range("G76").Value = '정품인증판매점'
$ws.Range("H76").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I76").Value = '삼성'
$ws.Range("J76").Value = ''
$ws.Range("M76").Value = '프린터'
$ws.Range("N76").Value = '레이저프린터'
$ws.Range("B80").Value = '에일리언테크놀로지아시아 뉴펜톡 3D펜 패키지 + PLA 필라멘트 5m 20색 세트'
$ws.Range("C80").Value = 'https://search.shopping.naver.com/gate.nhn?id=24173914522'
$ws.Range("D80").Value = 'https://shopping-phinf.pstatic.net/main_2417391/24173914522.20210216172007.jpg'
$ws.Range("E80").NumberFormat = "@"
$ws.Range("E80").Value = '102300'
$ws.Range("I80").Value = '펜톡'
$ws.Range("J80").Value = '에일리언테크놀로지아시아'
$ws.Range("N80").Value = '3D프린터'
$ws.Range("B81").Value = '캐논 LBP-6033 정품토너포함'
$ws.Range("C81").Value = 'https://search.shopping.naver.com/gate.nhn?id=7838416148'
$ws.Range("D81").Value = 'https://shopping-phinf.pstatic.net/main_7838416/7838416148.20141206171232.jpg'
$ws.Range("E81").NumberFormat = "@"
$ws.Range("E81").Value = '109900'
$ws.Range("G81").Value = '네이버'
$ws.Range("H81").Value = '일반 - 가격비교 상품'
$ws.Range("I81").Value = '캐논'
$ws.Range("J81").Value = '캐논'
$ws.Range("M81").Value = '프린터'
$ws.Range("N81").Value = '레이저프린터'
$ws.Range("B82").Value = '삼성 복합기렌탈 CLX-6260FR 사무용 컬러 레이저 프린터 임대 36개월'
$ws.Range("C82").Value = 'https://search.shopping.naver.com/gate.nhn?id=27492096087'
$ws.Range("D82").Value = 'https://shopping-phinf.pstatic.net/main_2749209/27492096087.2.jpg'
$ws.Range("E82").NumberFormat = "@"
$ws.Range("E82").Value = '40000'
$ws.Range("G82").Value = '오에이유'
$ws.Range("H82").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I82").Value = '삼성'
$ws.Range("J82").Value = '삼성전자'
$ws.Range("M82").Value = '복합기'
$ws.Range("N82").Value = '컬러레이저복합기'
